# Insert a new weekly record as row 25, pushing the existing rows 25-56
# down to 26-57 (matches how the source data feed prepends the latest
# observation to the top of the historical series).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("25:25").Insert()

$ws.Cells.Item(25, 1).Value = 11
$ws.Cells.Item(25, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(25, 3).Value = "Bíobío"
$ws.Cells.Item(25, 4).Value = 44967
$ws.Cells.Item(25, 5).Value = 8
$ws.Cells.Item(25, 6).Value = 100112030
$ws.Cells.Item(25, 7).Value = "Poroto granado"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 100
$ws.Cells.Item(25, 11).Value = 30000
$ws.Cells.Item(25, 12).Value = 32000
$ws.Cells.Item(25, 13).Value = 31000
$ws.Cells.Item(25, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(25, 15).Value = "Región Metropolitana"
$ws.Cells.Item(25, 16).Value = 1240
$ws.Cells.Item(25, 17).Value = 25
$ws.Cells.Item(25, 18).Value = "Hortaliza"

$ws.Cells.Item(25, 4).NumberFormat = $ws.Cells.Item(26, 4).NumberFormat
